$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados..." timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 13:06"

# --- Swap Malta / Crucero rows (keeps the sheet sorted descending by Casos totales) ---
# Row 155 held "Crucero", row 156 held "Malta"; after the update Malta's total (720)
# overtakes Crucero's (712), so their rows swap while keeping the list sorted.
$ws.Range("A155").Value = "Malta"
$ws.Range("A156").Value = "Crucero"

# --- Updated case numbers per country row ---
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4498475
$ws.Range("C4").Value = 132
$ws.Range("D4").Value = 2188954
$ws.Range("E4").Value = 2157178
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 152343

# Row 14
$ws.Range("B14").Value = 298909
$ws.Range("C14").Value = 2636
$ws.Range("D14").Value = 259116
$ws.Range("E14").Value = 23450
$ws.Range("G14").Value = 196
$ws.Range("H14").Value = 16343

# Row 43
$ws.Range("B43").Value = 59921
$ws.Range("C43").Value = 375
$ws.Range("D43").Value = 53202
$ws.Range("E43").Value = 6372

# Row 47
$ws.Range("B47").Value = 48235
$ws.Range("C47").Value = 1182
$ws.Range("D47").Value = 26446
$ws.Range("E47").Value = 19520
$ws.Range("G47").Value = 30
$ws.Range("H47").Value = 2269

# Row 55
$ws.Range("B55").Value = 34802
$ws.Range("C55").Value = 193
$ws.Range("E55").Value = 1823
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 1979

# Row 64
$ws.Range("E64").Value = 9776
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 128

# Row 67
$ws.Range("B67").Value = 19273
$ws.Range("C67").Value = 210
$ws.Range("D67").Value = 14021
$ws.Range("E67").Value = 5203

# Row 131
$ws.Range("B131").Value = 1861
$ws.Range("C131").Value = 4
$ws.Range("E131").Value = 28

# Row 155 (now Malta)
$ws.Range("B155").Value = 720
$ws.Range("C155").Value = 12
$ws.Range("D155").Value = 665
$ws.Range("E155").Value = 46
$ws.Range("H155").Value = 9

# Row 156 (now Crucero)
$ws.Range("B156").Value = 712
$ws.Range("D156").Value = 651
$ws.Range("E156").Value = 48
$ws.Range("H156").Value = 13

# Row 164
$ws.Range("B164").Value = 447
$ws.Range("C164").Value = 1
$ws.Range("E164").Value = 78
